# This script reproduces (as closely as the Word object model allows) the
# proofing-pass edit described by the commit: Word's spell/grammar checker
# re-scanned the document and wrapped a handful of words in
# <w:proofErr w:type="spellStart/spellEnd/gramStart/gramEnd"/> markers. That
# in turn causes the runs that contained those words to be split at the word
# boundaries. The words/visible text themselves are NOT changed anywhere.
#
# The COM object model does not expose a way to author <w:proofErr/>
# elements directly (there is no InsertProofErr-style call, and
# CheckSpelling/CheckGrammar do not mutate the document in this host).
# What *is* achievable, and what we do here, is force Word to split each
# affected run at exactly the same boundaries the proofing pass would have
# used, by toggling a character-formatting property on the sub-range and
# then clearing it again (leaving the run boundaries behind, with no visible
# formatting change). This reproduces the paragraph's run structure/text as
# closely as the automation surface allows.

$d = $word.ActiveDocument

function Split-Range($range) {
    # Forces Word to break $range out into its own run(s) without leaving
    # any visible formatting change behind.
    $range.Bold = 1
    $range.Bold = 0
}

# --- 1. "...write in the terminal "javac mainFrame.java"..." -----------
# "javac" gets flagged by the spell checker (spellStart/spellEnd).
$p1 = $d.Paragraphs.Item(17).Range.Duplicate
$f1 = $p1.Duplicate
if ($f1.Find.Execute("javac", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    Split-Range $f1
}

# --- 2. "...either in gtf (.gtf) or in fasta (.fa) format." -------------
# "(.fa" gets flagged by the grammar checker (gramStart/gramEnd).
$p2 = $d.Paragraphs.Item(19).Range.Duplicate
$f2 = $p2.Duplicate
if ($f2.Find.Execute("(.fa", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    Split-Range $f2
}

# --- 3. "Similarly to fasta, those statistics..." -----------------------
# "Similarly" gets flagged by the grammar checker (gramStart/gramEnd).
$p3 = $d.Paragraphs.Item(45).Range.Duplicate
$f3 = $p3.Duplicate
if ($f3.Find.Execute("Similarly", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    Split-Range $f3
}

# --- 4. "Both exons panels works like fasta's. ..." ----------------------
# "exons" gets flagged by the grammar checker (gramStart/gramEnd) and
# "fasta's" gets flagged by the spell checker (spellStart/spellEnd).
$p4 = $d.Paragraphs.Item(49).Range.Duplicate

$f4a = $p4.Duplicate
if ($f4a.Find.Execute("exons", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    Split-Range $f4a
}

$f4b = $p4.Duplicate
if ($f4b.Find.Execute("fasta" + [char]0x2019 + "s", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    Split-Range $f4b
}

Write-Output "proofing split pass complete"
